# Add a small "100" label textbox to Slide 5 (ppt/slides/slide5.xml),
# just after the two existing pictures, matching the target OOXML:
#
#   <p:sp>
#     <p:nvSpPr><p:cNvPr name="Google Shape;298;p16"/> ...txBox="1".../p:nvSpPr>
#     <p:spPr>
#       <a:xfrm><a:off x="6469600" y="1634125"/><a:ext cx="495600" cy="321600"/></a:xfrm>
#       <a:prstGeom prst="rect"/><a:noFill/><a:ln><a:noFill/></a:ln>
#     </p:spPr>
#     <p:txBody>
#       <a:bodyPr anchor="t" lIns="91425" rIns="91425" tIns="91425" bIns="91425" wrap="square"><a:noAutofit/></a:bodyPr>
#       <a:p>
#         <a:pPr lvl="0" algn="l"><a:spcBef><a:spcPts val="0"/></a:spcBef><a:spcAft><a:spcPts val="0"/></a:spcAft><a:buNone/></a:pPr>
#         <a:r><a:rPr lang="pt-BR" sz="1300"><a:solidFill><a:schemeClr val="dk2"/></a:solidFill>
#              <a:latin typeface="Nunito"/><a:ea typeface="Nunito"/><a:cs typeface="Nunito"/><a:sym typeface="Nunito"/></a:rPr>
#              <a:t>100</a:t></a:r>
#         <a:endParaRPr sz="1300">...same Nunito/dk2...</a:endParaRPr>
#       </a:p>
#     </p:txBody>
#   </p:sp>
#
# EMU -> point conversion (914400 EMU/in, 72 pt/in):
#   x  = 6469600 EMU = 509.41732283464563 pt
#   y  = 1634125 EMU = 128.67125984251967 pt
#   cx =  495600 EMU =  39.023622047244096 pt
#   cy =  321600 EMU =  25.322834645669293 pt

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

$shp = $s.Shapes.AddTextbox(1, 509.41732283464563, 128.67125984251967, 39.023622047244096, 25.322834645669293)
$shp.Name = "Google Shape;298;p16"

# <a:noFill/> shape fill, <a:ln><a:noFill/></a:ln> outline
$shp.Fill.Visible = 0
$shp.Line.Visible = 0

$tf = $shp.TextFrame
$tf.WordWrap = -1
$tf.AutoSize = 0
$tf.VerticalAnchor = 1
$tf.MarginLeft = 91425 / 12700.0
$tf.MarginRight = 91425 / 12700.0
$tf.MarginTop = 91425 / 12700.0
$tf.MarginBottom = 91425 / 12700.0

$tr = $tf.TextRange
$tr.Text = "100"
$tr.LanguageID = "pt-BR"
$tr.Font.Name = "Nunito"
$tr.Font.NameFarEast = "Nunito"
$tr.Font.NameComplexScript = "Nunito"
$tr.Font.Size = 13
$tr.Font.Color.ObjectThemeColor = 3  # msoThemeColorDark2 -> <a:schemeClr val="dk2"/>

$pf = $tr.ParagraphFormat
$pf.Alignment = 1            # ppAlignLeft -> algn="l"
$pf.Bullet.Visible = 0       # <a:buNone/>
$pf.SpaceBefore = 0          # <a:spcBef><a:spcPts val="0"/></a:spcBef>
$pf.SpaceAfter = 0           # <a:spcAft><a:spcPts val="0"/></a:spcAft>
$tr.IndentLevel = 1          # lvl="0" (1-based IndentLevel -> 0-based lvl)
